$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "93.282.47"
Set-TextValue $ws.Range("E2") "  +1.54%  "

Set-TextValue $ws.Range("D3") "3.122.52"
Set-TextValue $ws.Range("E3") "  -0.10%  "

Set-TextValue $ws.Range("E4") "  +0.01%  "

Set-TextValue $ws.Range("D5") "237.22"
Set-TextValue $ws.Range("E5") "  -3.34%  "

Set-TextValue $ws.Range("D6") "613.63"
Set-TextValue $ws.Range("E6") "  -0.80%  "

Set-TextValue $ws.Range("E7") "  -0.27%  "

Set-TextValue $ws.Range("D8") "0.391"
Set-TextValue $ws.Range("E8") "  +1.29%  "

Set-TextValue $ws.Range("E9") "  -0.06%  "

Set-TextValue $ws.Range("D10") "0.811"
Set-TextValue $ws.Range("E10") "  +8.06%  "

Set-TextValue $ws.Range("D11") "3.119.18"
Set-TextValue $ws.Range("E11") "  -0.06%  "

Set-TextValue $ws.Range("E12") "  -2.85%  "

Set-TextValue $ws.Range("D13") "0.0000245"
Set-TextValue $ws.Range("E13") "  -2.73%  "

Set-TextValue $ws.Range("D14") "34.86"
Set-TextValue $ws.Range("E14") "  -0.39%  "

Set-TextValue $ws.Range("D15") "92.988.03"
Set-TextValue $ws.Range("E15") "  +1.57%  "

Set-TextValue $ws.Range("D16") "5.44"
Set-TextValue $ws.Range("E16") "  -3.42%  "

Set-TextValue $ws.Range("D17") "3.709.16"
Set-TextValue $ws.Range("E17") "  +0.38%  "

Set-TextValue $ws.Range("D18") "3.129.74"
Set-TextValue $ws.Range("E18") "  +0.61%  "

Set-TextValue $ws.Range("D19") "3.77"
Set-TextValue $ws.Range("E19") "  +0.25%  "

Set-TextValue $ws.Range("D20") "14.63"
Set-TextValue $ws.Range("E20") "  -1.91%  "

Set-TextValue $ws.Range("D21") "5.95"
Set-TextValue $ws.Range("E21") "  +2.36%  "

Set-TextValue $ws.Range("D22") "0.0000203"
Set-TextValue $ws.Range("E22") "  -0.37%  "

Set-TextValue $ws.Range("D23") "441.77"
Set-TextValue $ws.Range("E23") "  -2.07%  "

Set-TextValue $ws.Range("D24") "9.11"
Set-TextValue $ws.Range("E24") "  -2.77%  "

Set-TextValue $ws.Range("B25") "NEARProtocol"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D25") "5.65"
Set-TextValue $ws.Range("E25") "  -3.97%  "

Set-TextValue $ws.Range("B26") "Aptos"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D26") "12.55"
Set-TextValue $ws.Range("E26") "  +6.84%  "

Set-TextValue $ws.Range("B27") "Litecoin"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D27") "85.93"
Set-TextValue $ws.Range("E27") "  -4.24%  "

Set-TextValue $ws.Range("B28") "WrappedeETH"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D28") "3.293.27"
Set-TextValue $ws.Range("E28") "  +0.08%  "

Set-TextValue $ws.Range("E29") "  +0.19%  "

Set-TextValue $ws.Range("D30") "0.179"
Set-TextValue $ws.Range("E30") "  +6.81%  "

Set-TextValue $ws.Range("D31") "0.232"
Set-TextValue $ws.Range("E31") "  +1.45%  "

Set-TextValue $ws.Range("E32") "  -13.28%  "

Set-TextValue $ws.Range("B33") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D33") "1.03"
Set-TextValue $ws.Range("E33") "  -1.72%  "

Set-TextValue $ws.Range("B34") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "9.18"
Set-TextValue $ws.Range("E34") "  -2.22%  "

Set-TextValue $ws.Range("D35") "8.17"
Set-TextValue $ws.Range("E35") "  +5.77%  "

Set-TextValue $ws.Range("E36") "  -9.14%  "

Set-TextValue $ws.Range("D37") "25.84"
Set-TextValue $ws.Range("E37") "  -1.92%  "

Set-TextValue $ws.Range("D38") "3.98"
Set-TextValue $ws.Range("E38") "  +0.74%  "

Set-TextValue $ws.Range("E39") "  -2.95%  "

Set-TextValue $ws.Range("E40") "  -0.85%  "

Set-TextValue $ws.Range("D42") "474.76"
Set-TextValue $ws.Range("E42") "  -3.27%  "

Set-TextValue $ws.Range("E43") "  -0.88%  "

Set-TextValue $ws.Range("D44") "3.32"
Set-TextValue $ws.Range("E44") "  -3.69%  "

Set-TextValue $ws.Range("E45") "  -0.03%  "

Set-TextValue $ws.Range("D46") "158.92"
Set-TextValue $ws.Range("E46") "  +0.63%  "

Set-TextValue $ws.Range("D47") "0.689"
Set-TextValue $ws.Range("E47") "  -1.89%  "

Set-TextValue $ws.Range("E48") "  -3.66%  "

Set-TextValue $ws.Range("E49") "  -1.59%  "

Set-TextValue $ws.Range("D50") "43.94"
Set-TextValue $ws.Range("E50") "  -0.58%  "

Set-TextValue $ws.Range("D51") "4.39"
Set-TextValue $ws.Range("E51") "  -1.20%  "
